$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B/C column values for rows 2-5 (ligand/receptor symbol swapped to Lgi3/Flot1)
$ws.Range("B2:B5").Value = "Lgi3"
$ws.Range("C2:C5").Value = "Flot1"

# Row 2 (FAPs -> ECs)
$ws.Range("D2").Value = "ECs"
$ws.Range("I2:J2").Value = 1
$ws.Range("M2").Value = 11.911367
$ws.Range("N2").Value = 35.734101
$ws.Range("O2:P2").Value = 0.2917564621783566
$ws.Range("Q2").Value = 18.293132008957
$ws.Range("R2").Value = 164.638188080613
$ws.Range("S2:T2").Value = 0.2917564621783566

# Row 3 (FAPs -> FAPs)
$ws.Range("I3:J3").Value = 1
$ws.Range("O3:P3").Value = 0.3312106030076451
$ws.Range("S3:T3").Value = 0.3312106030076451

# Row 4 (FAPs -> MuSCs)
$ws.Range("D4").Value = "MuSCs"
$ws.Range("I4:J4").Value = 1
$ws.Range("M4").Value = 7.266852666666666
$ws.Range("N4").Value = 21.800558
$ws.Range("O4:P4").Value = 0.1779939468910683
$ws.Range("Q4").Value = 11.16022158673933
$ws.Range("R4").Value = 100.441994280654
$ws.Range("S4:T4").Value = 0.1779939468910683

# Row 5 (FAPs -> Resolving-Mac)
$ws.Range("I5:J5").Value = 1
$ws.Range("M5").Value = 8.126046000000001
$ws.Range("N5").Value = 24.378138
$ws.Range("O5:P5").Value = 0.19903898792293
$ws.Range("Q5").Value = 12.479745791466
$ws.Range("R5").Value = 112.317712123194
$ws.Range("S5:T5").Value = 0.19903898792293

# Remove the old MuSCs-as-sender rows (previously rows 6-9)
$ws.Range("A6:A9").EntireRow.Delete()
